# Split the title run "AD and immigration" into three runs:
#   "AD" + " (American dream)" + " and immigration"
# so the document reads "AD (American dream) and immigration".

$d = $word.ActiveDocument

# Locate "AD" at the start of the title paragraph and collapse the
# found range to its end (i.e. the insertion point right after "AD").
$r = $d.Content
$r.Find.Execute("AD", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)

# Insert the new text right after "AD". InsertAfter repositions $r to
# cover exactly the text that was just inserted.
$r.InsertAfter(" (American dream)")

# Toggling a character property on/off forces the inserted text to stay
# in its own run (distinct from the neighbouring "AD" / " and immigration"
# runs) instead of being silently re-merged with them, while leaving the
# final formatting identical to the surrounding text.
$r.Bold = 1
$r.Bold = 0
